$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add two new bug rows (35 and 36) raised by "AD" on 2019-01-15 (serial 43480)
# Copy the date formatting from an existing "Date" cell (column B) so the
# new cells reuse the workbook's existing date style instead of creating a
# brand-new one.
$ws.Cells.Item(2, 2).Copy() | Out-Null
$ws.Cells.Item(35, 2).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(36, 2).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(35, 1).Value = "AD"
$ws.Cells.Item(35, 2).Value = 43480
$ws.Cells.Item(35, 3).Value = 'sumfig not works for ordered factor. Ordered factor class returns "ordered" "factor", fails in class checking. Suggestion: change variable.class == "factor" to any(variable.class == "factor")'

$ws.Cells.Item(36, 1).Value = "AD"
$ws.Cells.Item(36, 2).Value = 43480
$ws.Cells.Item(36, 3).Value = "SUGGESTION. sumby generates figure even set fig=F, it would be better sumby skip sumfig if fig set to F."

# Scroll / selection change to reflect reviewing the new rows
$ws.Application.ActiveWindow.ScrollRow = 25
$ws.Range("C39").Select()

# Set up the print page (adds pageSetup element with paper size / orientation)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
